$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Reword the "Low probability..." note in A3 (shared string text edit).
#    The original text contains a non-breaking space between "the" and "patient";
#    only the tail of the sentence changes.
$nbsp = [char]0x00A0
$newNote = "Low probability does NOT mean that the" + $nbsp + "patient does not have COVID-19! High probability means high probability for COVID-19."
$ws.Range("A3").Value = $newNote

# 2. Update the % formula in D10 so it shows blank until both smell inputs are filled in.
$ws.Range("D10").Formula = '=IF(AND(ISBLANK(D6),ISBLANK(D7)),"",CONCAT(ROUND(D9*100,0),"%"))'

# 3. Move the active selection to D10 (last selected cell in the saved workbook).
$null = $ws.Range("D10").Select()
